$wb = $excel.ActiveWorkbook

# Rename the "Eye tracking measures" sheet to "Eye_tracking_measures"
$eyeSheet = $wb.Worksheets.Item("Eye tracking measures")
$eyeSheet.Name = "Eye_tracking_measures"

# Trim trailing whitespace from specific cells on the "Studies" sheet
$studies = $wb.Worksheets.Item("Studies")

$cellsToTrim = @("C9", "C11", "N45", "N75", "F101", "J110", "N117", "N131", "J135", "J163", "F178", "F181", "N186", "F205", "N216", "F221", "N238")

foreach ($addr in $cellsToTrim) {
    $cell = $studies.Range($addr)
    $val = $cell.Value2
    if ($val -ne $null) {
        $cell.Value2 = $val.TrimEnd()
    }
}
